$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.NumberFormat = "General"
    $c.Style = "Normal"
}

Set-TextCell "D2" "35.058.27"
Set-TextCell "E2" "  +1.62%  "
Set-TextCell "D3" "1.857.32"
Set-TextCell "E3" "  +3.24%  "
Set-TextCell "E4" "  +0.16%  "
Set-TextCell "D5" "237.70"
Set-TextCell "E5" "  +3.77%  "
Set-TextCell "E6" "  +1.77%  "
Set-TextCell "E7" "  +0.06%  "
Set-TextCell "D8" "42.59"
Set-TextCell "E8" "  +8.90%  "
Set-TextCell "D9" "0.330"
Set-TextCell "E9" "  +3.62%  "
Set-TextCell "D10" "0.0695"
Set-TextCell "E10" "  +2.71%  "
Set-TextCell "E11" "  +0.36%  "
Set-TextCell "D12" "2.127.51"
Set-TextCell "E12" "  +3.32%  "
Set-TextCell "D13" "1.860.56"
Set-TextCell "E13" "  +3.41%  "
Set-TextCell "E14" "  +3.04%  "
Set-TextCell "E15" "  +3.29%  "
Set-TextCell "D16" "4.70"
Set-TextCell "E16" "  +3.62%  "
Set-TextCell "D17" "35.046.93"
Set-TextCell "E17" "  +2.10%  "
Set-TextCell "D18" "70.25"
Set-TextCell "E18" "  +1.88%  "
Set-TextCell "D19" "0.0₃0796"
Set-TextCell "E19" "  +2.41%  "
Set-TextCell "D20" "240.97"
Set-TextCell "E20" "  +0.71%  "
Set-TextCell "D21" "12.15"
Set-TextCell "E21" "  +3.41%  "
Set-TextCell "D22" "4.73"
Set-TextCell "E22" "  +1.42%  "
Set-TextCell "E23" "  -0.02%  "
Set-TextCell "E24" "  +1.90%  "
Set-TextCell "D25" "171.13"
Set-TextCell "E25" "  -1.14%  "
Set-TextCell "D26" "1.94"
Set-TextCell "E26" "  +30.83%  "
Set-TextCell "E27" "  +3.27%  "
Set-TextCell "D28" "17.71"
Set-TextCell "E28" "  +3.24%  "
Set-TextCell "D29" "0.124"
Set-TextCell "E29" "  +2.79%  "
Set-TextCell "B30" "BinanceUSD"
Set-TextCell "C30" "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextCell "D30" "1.01"
Set-TextCell "E30" "  +0.16%  "
Set-TextCell "B31" "Hedera"
Set-TextCell "C31" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell "D31" "0.0558"
Set-TextCell "E31" "  +3.32%  "
Set-TextCell "E32" "  +0.19%  "
Set-TextCell "E33" "  +3.28%  "
Set-TextCell "E34" "  +13.80%  "
Set-TextCell "D35" "1.63"
Set-TextCell "E35" "  +23.02%  "
Set-TextCell "D36" "0.787"
Set-TextCell "E36" "  +13.92%  "
Set-TextCell "D37" "1.30"
Set-TextCell "E37" "  +6.55%  "
Set-TextCell "E38" "  +12.92%  "
Set-TextCell "D39" "91.74"
Set-TextCell "E39" "  +1.39%  "
Set-TextCell "E40" "  +6.57%  "
Set-TextCell "D41" "1.352.90"
Set-TextCell "E41" "  +2.43%  "
Set-TextCell "D42" "14.86"
Set-TextCell "E42" "  +4.54%  "
Set-TextCell "D43" "2.33"
Set-TextCell "E43" "  +6.17%  "
Set-TextCell "D44" "12.87"
Set-TextCell "E44" "  +56.58%  "
Set-TextCell "E45" "  +0.21%  "
Set-TextCell "D46" "2.76"
Set-TextCell "E46" "  +2.02%  "
Set-TextCell "E47" "  +5.72%  "
Set-TextCell "D48" "6.44"
Set-TextCell "E48" "  +4.80%  "
Set-TextCell "D49" "2.040.58"
Set-TextCell "E49" "  +2.99%  "
Set-TextCell "E50" "  +3.44%  "
Set-TextCell "D51" "3.43"
Set-TextCell "E51" "  +18.45%  "
